$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# EA feedback pass on the SC-18 / mobile-code CCI cell (column G, "CCI").
#
# Before:
#   G2 and G3 both pointed at the same rich-text string:
#       "CCI-001695 ... SC-18 (3) ..."  (struck through, red)
#       "SC-18"                         (normal)
#
# After:
#   G2 is replaced outright with just the corrected CCI text.
#   G3 keeps the original text (now fully struck-through / red) and has the
#   corrected CCI appended right after it.
#
# NOTE: G2's value is written before G3's so that G2 reuses the original
# shared-string slot and G3's (longer, merged) text lands in a new slot --
# matching how the workbook's sharedStrings table is laid out after a real
# Excel edit/save cycle.
# ---------------------------------------------------------------------------

$oldCci = "CCI-001695`nThe information system prevents the execution of organization-defined unacceptable mobile code.`nNIST SP 800-53 :: SC-18 (3)`nNIST SP 800-53A :: SC-18 (3).1`nNIST SP 800-53 Revision 4 :: SC-18 (3)`n"
$newCciTitle = "CCI-001162`nThe organization establishes implementation guidance for acceptable mobile code and mobile code technologies."
$newCciRefs = "NIST SP 800-53::SC-18`nNIST SP 800-53A::SC-18`nNIST SP 800-53 Revision 4::SC-18"

# --- G2: replace with the corrected CCI text ---------------------------------
$g2Run1Text = $newCciTitle
$g2Run2Text = "`n"
$g2Run3Text = $newCciRefs + "`n"
$g2Text = $g2Run1Text + $g2Run2Text + $g2Run3Text

$ws.Range("G2").Value = $g2Text

$g2Len1 = $g2Run1Text.Length
$g2Len2 = $g2Run2Text.Length
$g2Len3 = $g2Run3Text.Length

$g2Run1 = $ws.Range("G2").Characters(1, $g2Len1)
$g2Run1.Font.Name = "Calibri (Body)"
$g2Run1.Font.Size = 12
$g2Run1.Font.Color = 255
$g2Run1.Font.Strikethrough = $false

$g2Run2 = $ws.Range("G2").Characters($g2Len1 + 1, $g2Len2)
$g2Run2.Font.Name = "Calibri"
$g2Run2.Font.Size = 12
$g2Run2.Font.Color = 255
$g2Run2.Font.Strikethrough = $true

$g2Run3 = $ws.Range("G2").Characters($g2Len1 + $g2Len2 + 1, $g2Len3)
$g2Run3.Font.Name = "Calibri"
$g2Run3.Font.Size = 12
$g2Run3.Font.Color = 255
$g2Run3.Font.Strikethrough = $false

# --- G3: append the corrected CCI block after the existing (now struck) text
$g3Text = $oldCci + $newCciTitle + "`n" + $newCciRefs
$ws.Range("G3").Value = $g3Text

$g3Len1 = $oldCci.Length
$g3Len2 = $g3Text.Length - $g3Len1

$g3Run1 = $ws.Range("G3").Characters(1, $g3Len1)
$g3Run1.Font.Name = "Calibri"
$g3Run1.Font.Size = 12
$g3Run1.Font.Color = 255
$g3Run1.Font.Strikethrough = $true

$g3Run2 = $ws.Range("G3").Characters($g3Len1 + 1, $g3Len2)
$g3Run2.Font.Name = "Calibri (Body)"
$g3Run2.Font.Size = 12
$g3Run2.Font.Color = 255
$g3Run2.Font.Strikethrough = $false

# Make red the cell's own (base) font too -- matches the new "red font /
# yellow fill" style now applied to both CCI cells.
$ws.Range("G2:G3").Font.Color = 255

# Reflect the last selection/active-cell position saved with the workbook.
[void]$ws.Range("G3").Select()
